# Word COM-interop script implementing the "Swap in Double Array in C++ (WIN)"
# commit: adds a practice-note bullet + example under the V26 section, and
# adds a new V27 "swap matrix rows/columns" exercise (with its demo video),
# renumbering the trailing placeholder paragraph from V27 -> V28.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1) Locate the paragraph that holds the "https://youtu.be/wjShmQlLVFE" link
#    (the end of the V26 block) and append the new bullet + "VD:" example
#    paragraphs right after it.
# ---------------------------------------------------------------------------
$v26LinkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*wjShmQlLVFE*") {
        $v26LinkPara = $p
    }
}
if ($v26LinkPara -eq $null) {
    throw "Could not find the V26 hyperlink paragraph (wjShmQlLVFE)"
}

$insertAfterV26 = $d.Range($v26LinkPara.Range.End - 1, $v26LinkPara.Range.End - 1)

$bulletAndVdXml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Must declare MAX_COL when call double Array in a function. </w:t></w:r></w:p>' + `
    '<w:p ' + $wNs + '><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>VD:</w:t></w:r><w:r><w:t xml:space="preserve"> void </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>input(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> a[][MAX_COL],</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> n, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> m)</w:t></w:r></w:p>'

$insertAfterV26.InsertXML($bulletAndVdXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Locate the trailing "V27: " paragraph (it currently carries the
#    _GoBack bookmark) and replace it with:
#       - V27: [Thực hành] Hoán vị hàng, cột của ma trận
#       - a new hyperlink paragraph (placeholder text, turned into a real
#         hyperlink below)
#       - a blank paragraph
#       - V28: <bookmark>
# ---------------------------------------------------------------------------
$v27Para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "V27:*") {
        $v27Para = $p
    }
}
if ($v27Para -eq $null) {
    throw "Could not find the V27 paragraph"
}

$placeholder = "zzzPLACEHOLDERzzz"

$v27ReplacementXml = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">V27: </w:t></w:r><w:r><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Thực</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hành</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">] </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Hoán</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>vị</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hàng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cột</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>của</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ma </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trận</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' + `
    '<w:p ' + $wNs + '><w:r><w:t>' + $placeholder + '</w:t></w:r></w:p>' + `
    '<w:p ' + $wNs + '/>' + `
    '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">V28: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$v27Para.Range.InsertXML($v27ReplacementXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Turn the placeholder paragraph's text into a real hyperlink pointing at
#    the new demo video, matching the Hyperlink character style used by the
#    other links in this document.
# ---------------------------------------------------------------------------
$linkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*$placeholder*") {
        $linkPara = $p
    }
}
if ($linkPara -eq $null) {
    throw "Could not find the placeholder paragraph for the new hyperlink"
}

$linkRange = $linkPara.Range
$linkRange.MoveEnd(1, -1) | Out-Null
$linkRange.Text = "https://youtu.be/_Y4TpGRxBdA"
$d.Hyperlinks.Add($linkRange, "https://youtu.be/_Y4TpGRxBdA") | Out-Null

Write-Output "edit applied"
